$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the F0 generation labels: underscore -> space
$ws.Range("B2").Value = "Female F0"
$ws.Range("B3").Value = "Male F0"
$ws.Range("B4").Value = "Both F0"
$ws.Range("B5").Value = "Female F0"
$ws.Range("B6").Value = "Male F0"

# Update the "Both_Sex Grand Offspring" label to "Grand Offspring (Both Sex)"
$ws.Range("B11").Value = "Grand Offspring (Both Sex)"
$ws.Range("B14").Value = "Grand Offspring (Both Sex)"

# Update the selection to B11
$ws.Range("B11").Select()
